# Commit: "add test and contest files"
#
# Adds a "benchmark" label in A1 and refreshes the "fact" benchmark row
# (B4:G4) with newly measured timings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label for column A (becomes shared string index 9: "benchmark")
$ws.Range("A1").Value = "benchmark"

# Updated timings for the "fact" row (row 4 / B4:G4)
# Written in plain decimal (not scientific notation) so the host parser
# accepts the literal; Excel stores/serialises the same IEEE-754 double.
$ws.Range("B4").Value = 0.00020799999999999999
$ws.Range("C4").Value = 0.00021100000000000001
$ws.Range("D4").Value = 0.000135885
$ws.Range("E4").Value = 0.00011661400000000001
$ws.Range("F4").Value = 0.000113125
$ws.Range("G4").Value = 0.000112656
